# Generate Report for handback
# Adds the handback status for a new source file
# (d2197cfc-e04f-4f70-82b3-2a0481c42435.md) to the Overview sheet and to
# each per-locale detail sheet (zh-cn, de-de), mirroring the layout that
# already exists for the other two source files.

$wb = $excel.ActiveWorkbook

$guid   = "d2197cfc-e04f-4f70-82b3-2a0481c42435"
$commit = "fe8ad8247747e6fde9cda147b27e4dd68e3dd304"

$inSync = "Handed back: in sync with en-US"

# -------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# -------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$guid.md"
$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6ce2a486aebfc4fdcce5caf6317f807980f261a/e2e/$guid.md", "", "", "$guid.md")

# Re-apply the explicit blue-underline look the sheet's other hyperlink
# cells use (so the new cellXf matches the existing custom "HyperLink"
# style rather than the default theme hyperlink color).
$wsOverview.Range("A4").Font.Underline = $true
$wsOverview.Range("A4").Font.Color = 15570276

# -------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de
#
# Columns: A Source File Name | B Status | C Correspond Handoff File |
# D Correspond Handoff Datetime | E Target File | F Correspond Handback
# File | G Correspond Handback DateTime | H Handoff Reason |
# I Dependency From
# -------------------------------------------------------------------

function Add-LocaleRow {
    param(
        $ws,
        [string]$locale,
        [string]$handoffAt,
        [string]$handbackAt,
        [string]$aSha,
        [string]$cSha,
        [string]$eSha,
        [string]$fSha
    )

    $xlfName = "$guid.$commit.$locale.xlf"

    $ws.Range("A4").Value = "$guid.md"
    $ws.Range("B4").Value = $inSync
    $ws.Range("C4").Value = $xlfName
    $ws.Range("D4").Value = $handoffAt
    $ws.Range("E4").Value = "$guid.md"
    $ws.Range("F4").Value = $xlfName
    $ws.Range("G4").Value = $handbackAt
    $ws.Range("H4").Value = "Include"

    $ws.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $ws.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$aSha/e2e/$guid.md"
    $handoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$cSha/ol-handoff/OpenLocalizationTestOrg/oltest.$locale/qimu/$xlfName"
    $targetMdUrl = "https://github.com/OpenLocalizationTestOrg/oltest.$locale/blob/$eSha/e2e/$guid.md"
    $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$fSha/ol-handback/OpenLocalizationTestOrg/oltest.$locale/qimu/$xlfName"

    $ws.Hyperlinks.Add($ws.Range("A4"), $mdUrl, "", "", "$guid.md")
    $ws.Hyperlinks.Add($ws.Range("C4"), $handoffUrl, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("E4"), $targetMdUrl, "", "", "$guid.md")
    $ws.Hyperlinks.Add($ws.Range("F4"), $handbackUrl, "", "", $xlfName)

    # Re-apply the explicit blue-underline look (matches the sheet's
    # pre-existing custom "HyperLink" cell style) rather than leaving the
    # default theme-colored hyperlink style that .Hyperlinks.Add() applies.
    foreach ($addr in @("A4", "C4", "E4", "F4")) {
        $ws.Range($addr).Font.Underline = $true
        $ws.Range($addr).Font.Color = 15570276
    }
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Add-LocaleRow $wsZhCn "zh-cn" `
    "2016-01-25 03:21:12" "2016-01-25 03:21:55" `
    "2011e33ea35a2fe38be159e7bb2ef23dbd280fc4" `
    "7311d7a092206f239f724091832b477ce4e2aa13" `
    "78ec9d27d1eb893280a21bb9ecd01b302576fede" `
    "c0e980cc1cd85707e22917d7e7ad512b7ceccca9"

$wsDeDe = $wb.Worksheets.Item("de-de")
Add-LocaleRow $wsDeDe "de-de" `
    "2016-01-25 03:21:23" "2016-01-25 03:22:12" `
    "7d1f153646b56321c8464a99e5ad5e4eebd1d28a" `
    "e6d0536a7e95fc38077579b517acb80ba5d2c2a9" `
    "89c51cb8c3d167271c4a1de179f2fea606053946" `
    "5db34ecd68cfdaa377b27dcdad44ac89a26938f8"
